$d = $word.ActiveDocument

# Locate the exact run of text that makes up the pitch paragraph's body
# (everything after "PITCH : ") without disturbing the rest of the
# paragraph (its pPr / rsid attributes must survive untouched).
$old = "PITCH : John Doe, un homme lambda dans un monde lambda. Venez briser cette monotonie en parsemant sa vie de mensonges...mais a quel prix ?"
$new = "Jean Dupont, un homme banal dans un monde banal, brise cette monotonie en parsemant sa vie de mensonges... Mais a quel prix ?"

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the pitch paragraph text to replace."
}

# $rng now spans the whole sentence. Re-anchor a range over just that
# span (excluding the paragraph mark) so InsertXML only rewrites the
# runs inside the paragraph, leaving <w:pPr> / paragraph rsids intact.
$target = $d.Range($rng.Start, $rng.End)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">PITCH : </w:t></w:r><w:r><w:t>' + $new + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
